$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")
$ws.Select()
$ws.Range("B5").Value = "Client wireframe to test server"
$ws.Range("B6").Value = "Create Local Server Scene"
